# "Source.xlsx" writer-output fixture gets regenerated: the "str"/"int"/
# "float" sample columns (old C:E) are dropped and the "decimal_2" column
# (old F) slides left to become the new column C. The sheet is renamed
# Source -> Data, column A is widened now that the sheet is narrower, and
# the saved selection/used-range shrink accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old str/int/float columns; deleting whole columns shifts the
# decimal_2 column (old F, with its 0.00_ number format) left into C.
$ws.Columns("C:E").Delete()

# Rename the worksheet.
$ws.Name = "Data"

# Column A gets wider now that the sheet only has 3 narrow columns.
$ws.Columns.Item(1).ColumnWidth = 23

# Selection moves to B19 in the saved view.
$ws.Range("B19").Select()
